$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (former Resolving-Mac target rows 14-17 plus EC-target rows already handled via data rewrite)
$ws.Rows("14:17").Delete()

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Il18"
$ws.Range("C2").Value2 = "Il18rap"
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.2184343333333333
$ws.Range("H2").Value2 = 0.655303
$ws.Range("I2").Value2 = 0.008416673064019609
$ws.Range("J2").Value2 = 0.00841667306401961
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.03475766666666667
$ws.Range("N2").Value2 = 0.104273
$ws.Range("O2").Value2 = 0.05507738448317117
$ws.Range("P2").Value2 = 0.05507738448317117
$ws.Range("Q2").Value2 = 0.007592267746555555
$ws.Range("R2").Value2 = 0.06833040971899999
$ws.Range("S2").Value2 = 0.0004635683384161584
$ws.Range("T2").Value2 = 0.0004635683384161584

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Il18"
$ws.Range("C3").Value2 = "Il18rap"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.2184343333333333
$ws.Range("H3").Value2 = 0.655303
$ws.Range("I3").Value2 = 0.008416673064019609
$ws.Range("J3").Value2 = 0.00841667306401961
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.2728686666666666
$ws.Range("N3").Value2 = 0.8186059999999999
$ws.Range("O3").Value2 = 0.4323907186158528
$ws.Range("P3").Value2 = 0.4323907186158528
$ws.Range("Q3").Value2 = 0.05960388529088888
$ws.Range("R3").Value2 = 0.5364349676179999
$ws.Range("S3").Value2 = 0.00363929131450613
$ws.Range("T3").Value2 = 0.003639291314506131

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Il18"
$ws.Range("C4").Value2 = "Il18rap"
$ws.Range("D4").Value2 = "Resolving-Mac"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.2184343333333333
$ws.Range("H4").Value2 = 0.655303
$ws.Range("I4").Value2 = 0.008416673064019609
$ws.Range("J4").Value2 = 0.00841667306401961
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.3234433333333333
$ws.Range("N4").Value2 = 0.9703299999999999
$ws.Range("O4").Value2 = 0.512531896900976
$ws.Range("P4").Value2 = 0.512531896900976
$ws.Range("Q4").Value2 = 0.07065112888777776
$ws.Range("R4").Value2 = 0.6358601599899999
$ws.Range("S4").Value2 = 0.00431381341109732
$ws.Range("T4").Value2 = 0.004313813411097321

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Il18"
$ws.Range("C5").Value2 = "Il18rap"
$ws.Range("D5").Value2 = "FAPs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 4.199828
$ws.Range("H5").Value2 = 12.599484
$ws.Range("I5").Value2 = 0.1618270290283213
$ws.Range("J5").Value2 = 0.1618270290283213
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.03475766666666667
$ws.Range("N5").Value2 = 0.104273
$ws.Range("O5").Value2 = 0.05507738448317117
$ws.Range("P5").Value2 = 0.05507738448317117
$ws.Range("Q5").Value2 = 0.1459762216813333
$ws.Range("R5").Value2 = 1.313785995132
$ws.Range("S5").Value2 = 0.008913009497562154
$ws.Range("T5").Value2 = 0.008913009497562156

# Row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Il18"
$ws.Range("C6").Value2 = "Il18rap"
$ws.Range("D6").Value2 = "MuSCs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 4.199828
$ws.Range("H6").Value2 = 12.599484
$ws.Range("I6").Value2 = 0.1618270290283213
$ws.Range("J6").Value2 = 0.1618270290283213
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.2728686666666666
$ws.Range("N6").Value2 = 0.8186059999999999
$ws.Range("O6").Value2 = 0.4323907186158528
$ws.Range("P6").Value2 = 0.4323907186158528
$ws.Range("Q6").Value2 = 1.146001466589333
$ws.Range("R6").Value2 = 10.314013199304
$ws.Range("S6").Value2 = 0.06997250537302431
$ws.Range("T6").Value2 = 0.06997250537302432

# Row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Il18"
$ws.Range("C7").Value2 = "Il18rap"
$ws.Range("D7").Value2 = "Resolving-Mac"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 4.199828
$ws.Range("H7").Value2 = 12.599484
$ws.Range("I7").Value2 = 0.1618270290283213
$ws.Range("J7").Value2 = 0.1618270290283213
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.3234433333333333
$ws.Range("N7").Value2 = 0.9703299999999999
$ws.Range("O7").Value2 = 0.512531896900976
$ws.Range("P7").Value2 = 0.512531896900976
$ws.Range("Q7").Value2 = 1.358406367746667
$ws.Range("R7").Value2 = 12.22565730972
$ws.Range("S7").Value2 = 0.08294151415773482
$ws.Range("T7").Value2 = 0.08294151415773483

# Row 8
$ws.Range("A8").Value2 = "MuSCs"
$ws.Range("B8").Value2 = "Il18"
$ws.Range("C8").Value2 = "Il18rap"
$ws.Range("D8").Value2 = "FAPs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 3.307112333333333
$ws.Range("H8").Value2 = 9.921336999999999
$ws.Range("I8").Value2 = 0.1274290669918512
$ws.Range("J8").Value2 = 0.1274290669918513
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.03475766666666667
$ws.Range("N8").Value2 = 0.104273
$ws.Range("O8").Value2 = 0.05507738448317117
$ws.Range("P8").Value2 = 0.05507738448317117
$ws.Range("Q8").Value2 = 0.1149475081112222
$ws.Range("R8").Value2 = 1.034527573001
$ws.Range("S8").Value2 = 0.007018459717041968
$ws.Range("T8").Value2 = 0.007018459717041969

# Row 9
$ws.Range("A9").Value2 = "MuSCs"
$ws.Range("B9").Value2 = "Il18"
$ws.Range("C9").Value2 = "Il18rap"
$ws.Range("D9").Value2 = "MuSCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 3.307112333333333
$ws.Range("H9").Value2 = 9.921336999999999
$ws.Range("I9").Value2 = 0.1274290669918512
$ws.Range("J9").Value2 = 0.1274290669918513
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.2728686666666666
$ws.Range("N9").Value2 = 0.8186059999999999
$ws.Range("O9").Value2 = 0.4323907186158528
$ws.Range("P9").Value2 = 0.4323907186158528
$ws.Range("Q9").Value2 = 0.9024073329135553
$ws.Range("R9").Value2 = 8.121665996221999
$ws.Range("S9").Value2 = 0.05509914584915421
$ws.Range("T9").Value2 = 0.05509914584915421

# Row 10
$ws.Range("A10").Value2 = "MuSCs"
$ws.Range("B10").Value2 = "Il18"
$ws.Range("C10").Value2 = "Il18rap"
$ws.Range("D10").Value2 = "Resolving-Mac"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 3.307112333333333
$ws.Range("H10").Value2 = 9.921336999999999
$ws.Range("I10").Value2 = 0.1274290669918512
$ws.Range("J10").Value2 = 0.1274290669918513
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.3234433333333333
$ws.Range("N10").Value2 = 0.9703299999999999
$ws.Range("O10").Value2 = 0.512531896900976
$ws.Range("P10").Value2 = 0.512531896900976
$ws.Range("Q10").Value2 = 1.069663436801111
$ws.Range("R10").Value2 = 9.626970931209998
$ws.Range("S10").Value2 = 0.06531146142565507
$ws.Range("T10").Value2 = 0.06531146142565508

# Row 11
$ws.Range("A11").Value2 = "Resolving-Mac"
$ws.Range("B11").Value2 = "Il18"
$ws.Range("C11").Value2 = "Il18rap"
$ws.Range("D11").Value2 = "FAPs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 18.22719966666667
$ws.Range("H11").Value2 = 54.681599
$ws.Range("I11").Value2 = 0.7023272309158078
$ws.Range("J11").Value2 = 0.7023272309158078
$ws.Range("K11").Value2 = 1
$ws.Range("L11").Value2 = 0.3333333333333333
$ws.Range("M11").Value2 = 0.03475766666666667
$ws.Range("N11").Value2 = 0.104273
$ws.Range("O11").Value2 = 0.05507738448317117
$ws.Range("P11").Value2 = 0.05507738448317117
$ws.Range("Q11").Value2 = 0.6335349302807778
$ws.Range("R11").Value2 = 5.701814372527
$ws.Range("S11").Value2 = 0.03868234693015089
$ws.Range("T11").Value2 = 0.03868234693015089

# Row 12
$ws.Range("A12").Value2 = "Resolving-Mac"
$ws.Range("B12").Value2 = "Il18"
$ws.Range("C12").Value2 = "Il18rap"
$ws.Range("D12").Value2 = "MuSCs"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 18.22719966666667
$ws.Range("H12").Value2 = 54.681599
$ws.Range("I12").Value2 = 0.7023272309158078
$ws.Range("J12").Value2 = 0.7023272309158078
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 0.2728686666666666
$ws.Range("N12").Value2 = 0.8186059999999999
$ws.Range("O12").Value2 = 0.4323907186158528
$ws.Range("P12").Value2 = 0.4323907186158528
$ws.Range("Q12").Value2 = 4.973631670110445
$ws.Range("R12").Value2 = 44.76268503099399
$ws.Range("S12").Value2 = 0.3036797760791681
$ws.Range("T12").Value2 = 0.3036797760791681

# Row 13
$ws.Range("A13").Value2 = "Resolving-Mac"
$ws.Range("B13").Value2 = "Il18"
$ws.Range("C13").Value2 = "Il18rap"
$ws.Range("D13").Value2 = "Resolving-Mac"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 18.22719966666667
$ws.Range("H13").Value2 = 54.681599
$ws.Range("I13").Value2 = 0.7023272309158078
$ws.Range("J13").Value2 = 0.7023272309158078
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 0.3234433333333333
$ws.Range("N13").Value2 = 0.9703299999999999
$ws.Range("O13").Value2 = 0.512531896900976
$ws.Range("P13").Value2 = 0.512531896900976
$ws.Range("Q13").Value2 = 5.895466217518889
$ws.Range("R13").Value2 = 53.05919595766999
$ws.Range("S13").Value2 = 0.3599651079064888
$ws.Range("T13").Value2 = 0.3599651079064888

